# Update the "4 cartes avec 15 valeur" (4card 15value) work item and its
# neighboring sub-version labels, marking it as tested/OK with completion
# dates, per commit message: "update 4card 15value for further development".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sub-version labels to include short descriptive suffixes.
$ws.Range("E8").Value  = "v2.1-3c15"
$ws.Range("E9").Value  = "v2.2-bre"
$ws.Range("E11").Value = "v3.1-4c15"
$ws.Range("E12").Value = "v3.2-carre"
$ws.Range("E13").Value = "v3.3-dp"

# Row 11 ("4 cartes avec 15 valeur") is now completed/tested: fill in the
# start/end dates (matching the existing date formatting used by the other
# rows) and mark the test result as OK.
$ws.Range("H9").Copy()
$ws.Range("H11:I11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H11").Value = 43167
$ws.Range("I11").Value = 43167
$ws.Range("J11").Value = "OK"

# Move the active selection to E13, matching where the author ended up.
$ws.Range("E13").Select()
